# Actualización 10 de Mayo
# Updates the statistics tables (Blancos, Reprobados, Aprobados, Por_Apro, Promedio)
# on sheets "Estadisticos 1P", "Estadisticos 2P" and "Estadisticos Final".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Estadisticos 1P
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Estadisticos 1P")

$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 16
$ws.Range("G2").Value = 44.44
$ws.Range("H2").Value = 6.3

$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 44.83
$ws.Range("H3").Value = 6.1

$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 28
$ws.Range("G4").Value = 77.78
$ws.Range("H4").Value = 7.8

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 16
$ws.Range("G5").Value = 72.73
$ws.Range("H5").Value = 7.5

# ---------------------------------------------------------------------------
# Estadisticos 2P
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Estadisticos 2P")

$ws.Range("D2").Value = 29
$ws.Range("E2").Value = 29
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 19.44
$ws.Range("H2").Value = 8.7

$ws.Range("D3").Value = 21
$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 27.59
$ws.Range("H3").Value = 7.4

$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 27
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = 8.3

$ws.Range("D5").Value = 12
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 45.45
$ws.Range("H5").Value = 8.8

# ---------------------------------------------------------------------------
# Estadisticos Final
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Estadisticos Final")

$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 16
$ws.Range("G2").Value = 44.44
$ws.Range("H2").Value = 6.3

$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 48.28
$ws.Range("H3").Value = 6.3

$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 30
$ws.Range("G4").Value = 83.33
$ws.Range("H4").Value = 8.2

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 16
$ws.Range("G5").Value = 72.73
$ws.Range("H5").Value = 7.6

Write-Host "Actualizacion aplicada"
